# Logged Week 15 and simulated Week 16
# Update row 2 ("H") values on both the OFF and DEF sheets with the
# latest Short/Deep Attempt & Completion counts.

$wb = $excel.ActiveWorkbook

$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 335
$wsOff.Range("C2").Value = 221
$wsOff.Range("D2").Value = 81
$wsOff.Range("E2").Value = 49

$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 418
$wsDef.Range("C2").Value = 283
$wsDef.Range("D2").Value = 88
$wsDef.Range("E2").Value = 44
